# Automatische test-sync: 2025-06-17 22:10:19
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# Append new row 48 to the Logs sheet
$logs.Range("A48").Value = "Afmelding nieuwsbrief"
$logs.Range("B48").Value = "mailmind.test@zohomail.eu"
$logs.Range("C48").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D48").Value = "Afmelding"
$logs.Range("F48").Value = "2025-06-17 22:09:47"
$logs.Range("G48").Value = "Nee"

# Extend the conditional formatting ranges to include the new row
$dRange = $logs.Range("D2:D48")
$dConditions = $logs.Range("D2").FormatConditions
for ($i = 1; $i -le $dConditions.Count; $i++) {
    $dConditions.Item($i).ModifyAppliesToRange($dRange)
}

$gRange = $logs.Range("G2:G48")
$gConditions = $logs.Range("G2").FormatConditions
for ($i = 1; $i -le $gConditions.Count; $i++) {
    $gConditions.Item($i).ModifyAppliesToRange($gRange)
}

# Update the Afmelding count on the Dashboard sheet
$dashboard.Range("B4").Value = 9
